# Auto-generated Excel COM-interop script
# Scheduled-runner refresh: updates Universalis market-price columns
# (currentAveragePrice / currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ -- columns H-N) for the rows whose item prices moved
# since the last sync. Values below are literal data snapshots, not
# formulas, matching how this sheet is populated upstream.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 117.13
$ws.Range("I15").Value = 117.13
$ws.Range("K15").Value = 351.39
$ws.Range("M15").Value = -182.39
$ws.Range("H62").Value = 3229.8
$ws.Range("I62").Value = 3233
$ws.Range("J62").Value = 3225
$ws.Range("K62").Value = 3233
$ws.Range("L62").Value = 3225
$ws.Range("M62").Value = -2609
$ws.Range("N62").Value = -4473
$ws.Range("H65").Value = 3229.8
$ws.Range("I65").Value = 3233
$ws.Range("J65").Value = 3225
$ws.Range("K65").Value = 16165
$ws.Range("L65").Value = 16125
$ws.Range("M65").Value = -13045
$ws.Range("N65").Value = -22365
$ws.Range("H74").Value = 3359.0908
$ws.Range("I74").Value = 2957.1428
$ws.Range("J74").Value = 4062.5
$ws.Range("K74").Value = 2957.1428
$ws.Range("L74").Value = 4062.5
$ws.Range("M74").Value = -2021.1428
$ws.Range("N74").Value = -5934.5
$ws.Range("H77").Value = 3359.0908
$ws.Range("I77").Value = 2957.1428
$ws.Range("J77").Value = 4062.5
$ws.Range("K77").Value = 14785.714
$ws.Range("L77").Value = 20312.5
$ws.Range("M77").Value = -10105.714
$ws.Range("N77").Value = -29672.5
$ws.Range("H137").Value = 1234.7826
$ws.Range("I137").Value = 1161.9524
$ws.Range("J137").Value = 1999.5
$ws.Range("K137").Value = 3485.857199999999
$ws.Range("L137").Value = 5998.5
$ws.Range("M137").Value = -935.8571999999995
$ws.Range("N137").Value = -11098.5
$ws.Range("H138").Value = 3410.1785
$ws.Range("J138").Value = 3501.8142
$ws.Range("L138").Value = 10505.4426
$ws.Range("N138").Value = -20785.4426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2738
$ws.Range("J63").Value = 3047.5
$ws.Range("L63").Value = 3047.5
$ws.Range("N63").Value = -4419.5
$ws.Range("H66").Value = 2738
$ws.Range("J66").Value = 3047.5
$ws.Range("L66").Value = 15237.5
$ws.Range("N66").Value = -22101.5
$ws.Range("H74").Value = 860.3333
$ws.Range("I74").Value = 496.8
$ws.Range("J74").Value = 1314.75
$ws.Range("K74").Value = 496.8
$ws.Range("L74").Value = 1314.75
$ws.Range("M74").Value = 377.2
$ws.Range("N74").Value = -3062.75
$ws.Range("H77").Value = 860.3333
$ws.Range("I77").Value = 496.8
$ws.Range("J77").Value = 1314.75
$ws.Range("K77").Value = 2484
$ws.Range("L77").Value = 6573.75
$ws.Range("M77").Value = 1884
$ws.Range("N77").Value = -15309.75
$ws.Range("H110").Value = 34556420
$ws.Range("I110").Value = 40085196
$ws.Range("J110").Value = 1562.5
$ws.Range("K110").Value = 40085196
$ws.Range("L110").Value = 1562.5
$ws.Range("M110").Value = -40083151
$ws.Range("N110").Value = -5652.5
$ws.Range("H122").Value = 3682
$ws.Range("I122").Value = 3203
$ws.Range("K122").Value = 9609
$ws.Range("M122").Value = -7159

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21061.715
$ws.Range("J82").Value = 34485
$ws.Range("L82").Value = 34485
$ws.Range("N82").Value = -35251
$ws.Range("H85").Value = 21061.715
$ws.Range("J85").Value = 34485
$ws.Range("L85").Value = 34485
$ws.Range("N85").Value = -37137
$ws.Range("H86").Value = 87830
$ws.Range("I86").Value = 103471.82
$ws.Range("K86").Value = 103471.82
$ws.Range("M86").Value = -102348.82
$ws.Range("H89").Value = 87830
$ws.Range("I89").Value = 103471.82
$ws.Range("K89").Value = 517359.1
$ws.Range("M89").Value = -511743.1
$ws.Range("H105").Value = 251964
$ws.Range("I105").Value = 202298
$ws.Range("J105").Value = 334740.66
$ws.Range("K105").Value = 202298
$ws.Range("L105").Value = 334740.66
$ws.Range("M105").Value = -200551
$ws.Range("N105").Value = -338234.66
$ws.Range("H132").Value = 63000
$ws.Range("J132").Value = 63000
$ws.Range("L132").Value = 63000
$ws.Range("N132").Value = -73120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22786.549
$ws.Range("I31").Value = 1136.775
$ws.Range("J31").Value = 50721.742
$ws.Range("K31").Value = 1136.775
$ws.Range("L31").Value = 50721.742
$ws.Range("M31").Value = -841.7750000000001
$ws.Range("N31").Value = -51311.742
$ws.Range("H34").Value = 22786.549
$ws.Range("I34").Value = 1136.775
$ws.Range("J34").Value = 50721.742
$ws.Range("K34").Value = 1136.775
$ws.Range("L34").Value = 50721.742
$ws.Range("M34").Value = -934.7750000000001
$ws.Range("N34").Value = -51125.742
$ws.Range("H51").Value = 7924.8335
$ws.Range("J51").Value = 7924.8335
$ws.Range("L51").Value = 7924.8335
$ws.Range("N51").Value = -9396.833500000001
$ws.Range("H61").Value = 7924.8335
$ws.Range("J61").Value = 7924.8335
$ws.Range("L61").Value = 7924.8335
$ws.Range("N61").Value = -8620.833500000001
$ws.Range("H99").Value = 6799.696
$ws.Range("I99").Value = 2040
$ws.Range("K99").Value = 2040
$ws.Range("M99").Value = -542
$ws.Range("H122").Value = 1400
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 4200
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9100
$ws.Range("H126").Value = 6799.696
$ws.Range("I126").Value = 2040
$ws.Range("K126").Value = 6120
$ws.Range("M126").Value = -3650
$ws.Range("H134").Value = 994.67566
$ws.Range("I134").Value = 543.1852
$ws.Range("J134").Value = 2213.7
$ws.Range("K134").Value = 1629.5556
$ws.Range("L134").Value = 6641.099999999999
$ws.Range("M134").Value = 905.4443999999999
$ws.Range("N134").Value = -11711.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 734434.2
$ws.Range("J131").Value = 817271.9399999999
$ws.Range("L131").Value = 2451815.82
$ws.Range("N131").Value = -2461895.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1354.6923
$ws.Range("I113").Value = 873
$ws.Range("J113").Value = 1916.6666
$ws.Range("K113").Value = 873
$ws.Range("L113").Value = 1916.6666
$ws.Range("M113").Value = 1297
$ws.Range("N113").Value = -6256.6666
$ws.Range("H122").Value = 4415.8
$ws.Range("I122").Value = 4879.857
$ws.Range("K122").Value = 14639.571
$ws.Range("M122").Value = -12189.571
$ws.Range("H126").Value = 3246.4
$ws.Range("I126").Value = 3186.0557
$ws.Range("J126").Value = 3401.5715
$ws.Range("K126").Value = 9558.167099999999
$ws.Range("L126").Value = 10204.7145
$ws.Range("M126").Value = -7088.167099999999
$ws.Range("N126").Value = -15144.7145
$ws.Range("H132").Value = 2611.3777
$ws.Range("I132").Value = 1832.2188
$ws.Range("J132").Value = 4529.3076
$ws.Range("K132").Value = 5496.6564
$ws.Range("L132").Value = 13587.9228
$ws.Range("M132").Value = -2966.6564
$ws.Range("N132").Value = -18647.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8012.8335
$ws.Range("J41").Value = 8615.4
$ws.Range("L41").Value = 8615.4
$ws.Range("N41").Value = -9395.4
$ws.Range("H137").Value = 48357.5
$ws.Range("J137").Value = 48357.5
$ws.Range("L137").Value = 48357.5
$ws.Range("N137").Value = -58557.5
